# Adds devaluation-related fields to the "DepositCalculationData" class map
# documented on Лист2 (predictions + aggregate "current"/"estimated" fields):
#   - DepositStates State            (was missing a listed row of its own)
#   - decimal TotalPercentInUsd
#   - decimal CurrentDevaluationInUsd
#   - decimal EstimatedCurrencyRateOnFinish
#   - decimal EstimatedDevaluationInUsd

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист2")

# Insert 4 new blank rows at the positions (final row numbers) where the new
# fields are introduced. Inserting top-to-bottom is safe because each Insert
# only pushes rows *below* the insertion point further down.
$newRowNumbers = @(23, 41, 45, 48)
foreach ($r in $newRowNumbers) {
    $ws.Rows.Item($r).Insert()
}

# Re-write the full class-member listing for rows 23-50 with the updated
# field set (existing fields shift down to make room for the new ones, so
# it's simplest to just restate the whole block).
$members = @(
    @(23, 3, "          DepositStates State"),
    @(24, 3, "          List<"),
    @(25, 4, "DateTime Date"),
    @(26, 4, "decimal Balance"),
    @(27, 4, "decimal DepoRate"),
    @(28, 4, "decimal DayProcents"),
    @(29, 4, "decimal NotPaidProcents"),
    @(30, 4, "decimal CurrencyRate"),
    @(31, 4, "decimal DayDevaluation"),
    @(32, 3, "          List<"),
    @(33, 4, "DateTime Timestamp"),
    @(34, 4, "TransactionType"),
    @(35, 4, "decimal Amount"),
    @(36, 4, "Currency"),
    @(37, 4, "decimal AmountInUsd"),
    @(38, 4, "string Comment"),
    @(39, 3, "          decimal TotalMyIns"),
    @(40, 3, "          decimal TotalPercent"),
    @(41, 3, "          decimal TotalPercentInUsd"),
    @(42, 3, "          decimal TotalMyOuts"),
    @(43, 3, "          decimal CurrentBalance"),
    @(44, 3, "          decimal CurrentProfitInUsd"),
    @(45, 3, "          decimal CurrentDevaluationInUsd"),
    @(46, 3, "          decimal EstimatedProcentsInThisMonth"),
    @(47, 3, "          decimal EstimatedProcents"),
    @(48, 3, "          decimal EstimatedCurrencyRateOnFinish"),
    @(49, 3, "          decimal EstimatedDevaluationInUsd"),
    @(50, 3, "          decimal EstimatedProfitInUsd")
)

foreach ($m in $members) {
    $row = $m[0]
    $col = $m[1]
    $text = $m[2]
    $cell = $ws.Cells.Item($row, $col)
    # The freshly-inserted rows pick up formatting from their neighbour; none
    # of these label rows carry explicit cell styling in the original sheet,
    # so strip whatever got inherited.
    $cell.ClearFormats()
    $cell.Value = $text
}

$ws.Range("F9").Select()
